# "Generate Report for Handback" - refresh the handback-status report.
#
# A re-run of the report generator only updated the timestamps for the
# "43bb357f-9e76-4b70-ac86-144c37b4199c" file's handoff/handback round-trip
# (row 2 on each per-locale sheet); the "bcaa9fb7-3425-4373-99c9-edffb332bcce"
# rows are untouched. Update:
#   - zh-cn!H2 (Correspond Handoff Datetime)   -> 2016-09-03 12:52:24
#   - zh-cn!K2 (Correspond Handback DateTime)  -> 2016-09-03 12:52:42
#   - de-de!H2 (Correspond Handoff Datetime)   -> 2016-09-03 12:52:28
#   - de-de!K2 (Correspond Handback DateTime)  -> 2016-09-03 12:52:49
#   - Overview!G2 (Latest HO Xliff Generate Date, mirrors de-de handoff time)
#                                               -> 2016-09-03 12:52:28

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-09-03 12:52:24"
$zhcn.Range("K2").Value = "2016-09-03 12:52:42"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-09-03 12:52:28"
$dede.Range("K2").Value = "2016-09-03 12:52:49"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-09-03 12:52:28"
